# CrewAI Robust Backend Ready!
# Re-labels the header row (lower-cased / renamed columns), shifts the
# numeric columns D/E one position left (D<-old E, E<-old F) and fills in
# a brand-new "climate change" metric in column F, then documents every
# header cell with a reviewer comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row rename
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "industry"
$ws.Cells.Item(1, 2).Value = "unit"
$ws.Cells.Item(1, 3).Value = "process"
$ws.Cells.Item(1, 4).Value = "carbon (kg CO2 eq)"
$ws.Cells.Item(1, 5).Value = "ced (MJ)"
$ws.Cells.Item(1, 6).Value = "climate change (kg CO2 eq)"
$ws.Cells.Item(1, 7).Value = "region"

# ---------------------------------------------------------------------
# 2. Data rows 2-27: column D gets the old column-E value, column E gets
#    the old column-F value, and column F gets the brand new
#    "climate change" figure.
# ---------------------------------------------------------------------
$data = @(
  @(2, [double]"0.40912154", [double]"3.5739869", [double]"1.1407394e-05"),
  @(3, [double]"0.9029999999999999", [double]"4.0130628", [double]"2.5178036e-05"),
  @(4, [double]"0.1221332466666667", [double]"1.3685633", [double]"3.4053989e-06"),
  @(5, [double]"0.4556034466666667", [double]"47.209177", [double]"1.2703433e-05"),
  @(6, [double]"1.731465466666667", [double]"21.29503", [double]"4.8277852e-05"),
  @(7, [double]"0.2776", [double]"4.03118", [double]"7.740224599999999e-06"),
  @(8, [double]"0.2232346266666667", [double]"3.2418132", [double]"6.2243737e-06"),
  @(9, [double]"0.361", [double]"5.242442", [double]"1.0065638e-05"),
  @(10, [double]"0.1335", [double]"1.938687", [double]"3.7223342e-06"),
  @(11, [double]"309.6184", [double]"1518.0958", [double]"0.008632982500000001"),
  @(12, [double]"450.4321333333334", [double]"3101.7854", [double]"0.012559243"),
  @(13, [double]"0.9279452666666668", [double]"36.340984", [double]"2.5873577e-05"),
  @(14, [double]"1.24", [double]"4.72", [double]"3.457449e-05"),
  @(15, [double]"1.2686458", [double]"1.034375", [double]"3.5373212e-05"),
  @(16, [double]"0.00297", [double]"0.034751155", [double]"8.2811481e-08"),
  @(17, [double]"0.8532838", [double]"5.0228714", [double]"2.3791816e-05"),
  @(18, [double]"0.01005431133333333", [double]"0.14197459", [double]"2.8034087e-07"),
  @(19, [double]"0.008451014", [double]"0.11135655", [double]"2.356367e-07"),
  @(20, [double]"0.000224564", [double]"0.0034906813", [double]"6.2614402e-09"),
  @(21, [double]"5.2387586", [double]"81.754907", [double]"0.00014607049"),
  @(22, [double]"1.88", [double]"36.3862", [double]"5.2419388e-05"),
  @(23, [double]"0.00442357", [double]"0.048757127", [double]"1.2334087e-07"),
  @(24, [double]"3.742187533333334", [double]"96.29328099999999", [double]"0.00010434212"),
  @(25, [double]"4.363466933333333", [double]"98.83900800000001", [double]"0.00012166504"),
  @(26, [double]"2.260869533333334", [double]"39.65113", [double]"6.3039042e-05"),
  @(27, [double]"3.204360266666667", [double]"52.532058", [double]"8.9346067e-05")
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 4).Value = $row[1]
  $ws.Cells.Item($r, 5).Value = $row[2]
  $ws.Cells.Item($r, 6).Value = $row[3]
}

# ---------------------------------------------------------------------
# 3. Reviewer comments on every header cell
# ---------------------------------------------------------------------
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
